# Apply updated cryptos list values (price + 1h volume change) per the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.377.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "'1.550.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'209.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "'0.480"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'23.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").Value = "'0.0582"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").Value = "'0.0888"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "'1.771.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").Value = "'1.565.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "'28.372.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").Value = "'0.509"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").Value = "'60.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("D18").Value = "'228.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").Value = "'0.0₃0675"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("D25").Value = "'151.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'14.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -3.34%  "
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("E31").Value = "  -4.96%  "
$ws.Range("D32").Value = "'3.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("D33").Value = "'1.384.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").Value = "'1.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("D37").Value = "'2.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").Value = "'0.511"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'0.773"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "'0.0456"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").Value = "'5.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("D47").Value = "'1.684.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").Value = "'0.871"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.31%  "
$ws.Range("D49").Value = "'85.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "'43.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.44%  "
$ws.Range("D51").Value = "'0.0₆0102"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
